$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Drop the EuropeMarket and AsiaMarket sheets, keep only USMarket ---
$wb.Worksheets.Item("EuropeMarket").Delete() | Out-Null
$wb.Worksheets.Item("AsiaMarket").Delete() | Out-Null

# --- Rename the remaining sheet to AMG ---
$ws = $wb.Worksheets.Item("USMarket")
$ws.Name = "AMG"

# --- Replace the row 2 / row 3 text, clear rows 4-7 (keep their style) ---
$ws.Range("A2").Value = "MERCEDES-AMG ENGINES"
$ws.Range("A3").Value = "PERFORMANCE 6/6"
$ws.Range("A4").Value = $null
$ws.Range("A5").Value = $null
$ws.Range("A6").Value = $null
$ws.Range("A7").Value = $null

# --- Column A width tweak (closest reachable value to the recorded bestFit width) ---
$ws.Columns("A").ColumnWidth = 21.8

# --- Selection / active sheet view ---
$ws.Activate()
$ws.Range("A4").Select() | Out-Null
